# Natmi following Dr Hou advice
#
# The Mcam-Mcam LR-pair sheet recomputed ligand/receptor-expressing cell counts
# (column E "Ligand-expressing cells" and column K "Receptor-expressing cells")
# from 1 to 3 for every data row, which in turn changes every downstream
# average/total expression, derived-specificity and edge-weight column
# (G, H, I, J, M, N, O, P, Q, R, S, T). Columns A-D (cluster/gene labels) and
# F, L (detection rates, unchanged at 1) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

$rowsData = @(
    @{ Row = 2; Vals = @(3,87.038094,261.114282,0.7371871251537216,0.7371871251537216,3,87.038094,261.114282,0.7371871251537216,0.7371871251537216,7575.629807152836,68180.66826437552,0.5434448574924089,0.5434448574924089) },
    @{ Row = 3; Vals = @(3,87.038094,261.114282,0.7371871251537216,0.7371871251537216,3,2.255465,6.766394999999999,0.01910312695076754,0.01910312695076754,196.31137468371,1766.80237215339,0.0140825792382829,0.01408257923828291) },
    @{ Row = 4; Vals = @(3,87.038094,261.114282,0.7371871251537216,0.7371871251537216,3,0.5269253333333334,1.580776,0.004462903009464643,0.004462903009464643,45.862576693648,412.763190242832,0.003289994639387133,0.003289994639387133) },
    @{ Row = 5; Vals = @(3,87.038094,261.114282,0.7371871251537216,0.7371871251537216,3,28.247359,84.74207700000001,0.2392468448860462,0.2392468448860462,2458.596287893746,22127.36659104372,0.1763696937836428,0.1763696937836428) },
    @{ Row = 6; Vals = @(3,2.255465,6.766394999999999,0.01910312695076754,0.01910312695076754,3,87.038094,261.114282,0.7371871251537216,0.7371871251537216,196.31137468371,1766.80237215339,0.0140825792382829,0.01408257923828291) },
    @{ Row = 7; Vals = @(3,2.255465,6.766394999999999,0.01910312695076754,0.01910312695076754,3,2.255465,6.766394999999999,0.01910312695076754,0.01910312695076754,5.087122366224998,45.78410129602499,0.0003649294592971411,0.0003649294592971412) },
    @{ Row = 8; Vals = @(3,2.255465,6.766394999999999,0.01910312695076754,0.01910312695076754,3,0.5269253333333334,1.580776,0.004462903009464643,0.004462903009464643,1.188461646946666,10.69615482252,0.00008525540275876558,0.0000852554027587656) },
    @{ Row = 9; Vals = @(3,2.255465,6.766394999999999,0.01910312695076754,0.01910312695076754,3,28.247359,84.74207700000001,0.2392468448860462,0.2392468448860462,63.71092956693499,573.398366102415,0.00457036285042873,0.004570362850428731) },
    @{ Row = 10; Vals = @(3,0.5269253333333334,1.580776,0.004462903009464643,0.004462903009464643,3,87.038094,261.114282,0.7371871251537216,0.7371871251537216,45.862576693648,412.763190242832,0.003289994639387133,0.003289994639387133) },
    @{ Row = 11; Vals = @(3,0.5269253333333334,1.580776,0.004462903009464643,0.004462903009464643,3,2.255465,6.766394999999999,0.01910312695076754,0.01910312695076754,1.188461646946666,10.69615482252,0.00008525540275876558,0.0000852554027587656) },
    @{ Row = 12; Vals = @(3,0.5269253333333334,1.580776,0.004462903009464643,0.004462903009464643,3,0.5269253333333334,1.580776,0.004462903009464643,0.004462903009464643,0.2776503069084444,2.498852762176,0.00001991750327188857,0.00001991750327188857) },
    @{ Row = 13; Vals = @(3,0.5269253333333334,1.580776,0.004462903009464643,0.004462903009464643,3,28.247359,84.74207700000001,0.2392468448860462,0.2392468448860462,14.88424905686134,133.958241511752,0.001067735464046856,0.001067735464046856) },
    @{ Row = 14; Vals = @(3,28.247359,84.74207700000001,0.2392468448860462,0.2392468448860462,3,87.038094,261.114282,0.7371871251537216,0.7371871251537216,2458.596287893746,22127.36659104372,0.1763696937836428,0.1763696937836428) },
    @{ Row = 15; Vals = @(3,28.247359,84.74207700000001,0.2392468448860462,0.2392468448860462,3,2.255465,6.766394999999999,0.01910312695076754,0.01910312695076754,63.71092956693499,573.398366102415,0.00457036285042873,0.004570362850428731) },
    @{ Row = 16; Vals = @(3,28.247359,84.74207700000001,0.2392468448860462,0.2392468448860462,3,0.5269253333333334,1.580776,0.004462903009464643,0.004462903009464643,14.88424905686134,133.958241511752,0.001067735464046856,0.001067735464046856) },
    @{ Row = 17; Vals = @(3,28.247359,84.74207700000001,0.2392468448860462,0.2392468448860462,3,28.247359,84.74207700000001,0.2392468448860462,0.2392468448860462,797.9132904748811,7181.21961427393,0.05723905278792783,0.05723905278792785) }
)

foreach ($rd in $rowsData) {
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $cellRef = "$($colLetters[$i])$($rd.Row)"
        $ws.Range($cellRef).Value = $rd.Vals[$i]
    }
}